{"js": "// The document contains \"tcn\" (transcription) markers where the XML-like\n// tag text, e.g. \"<id>p112r_4</id>\", was previously split across three\n// runs: \"<id>\", \"p112r_4\", \"</id>\". This edit merges each split marker\n// back into a single run (keeping the opening run's formatting), for the\n// two newly downloaded ids \"p112r_4\" and \"p112v_1\" \u2014 matching the commit\n// \"add newly downloaded tc, tcn. tl\".\nconst ids = [\"p112r_4\", \"p112v_1\"];\n\nfor (const id of ids) {\n  const fullText = \"<id>\" + id + \"</id>\";\n  const results = context.document.body.search(fullText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const hit of results.items) {\n    // Replacing the whole matched range with the identical text collapses\n    // the (up to three) runs spanned by the match into a single run that\n    // carries the formatting of the first run in the match.\n    hit.insertText(fullText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# The document contains \"tcn\" (transcription) markers where the XML-like\n# tag text, e.g. \"<id>p112r_4</id>\", was previously split across three\n# runs: \"<id>\", \"p112r_4\", \"</id>\". This edit merges each split marker\n# back into a single run (keeping the opening run's formatting), for the\n# two newly downloaded ids \"p112r_4\" and \"p112v_1\" \u2014 matching the commit\n# \"add newly downloaded tc, tcn. tl\".\n\n$d = $word.ActiveDocument\n$ids = @(\"p112r_4\", \"p112v_1\")\n\nforeach ($id in $ids) {\n    $target = \"<id>\" + $id + \"</id>\"\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $target\n    $find.Replacement.Text = $target\n\n    # Replacing the whole matched text with itself collapses the (up to\n    # three) runs spanned by the match into a single run carrying the\n    # formatting of the match's first run.\n    $find.Execute(\n        $target,        # FindText\n        $false,         # MatchCase\n        $false,         # MatchWholeWord\n        $false,         # MatchWildcards\n        $false,         # MatchSoundsLike\n        $false,         # MatchAllWordForms\n        $true,          # Forward\n        1,              # Wrap (wdFindContinue)\n        $false,         # Format\n        $target,        # ReplaceWith\n        2               # Replace (wdReplaceAll)\n    ) | Out-Null\n}\n"}
